# "updated test for add user"
# Adds a "Read User Table" section header above the existing user-button /
# read-table steps, and appends a new "Add user" section with its own
# keyword rows, mirroring the layout/style already used in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Make room for the new "Read User Table" section header: push the
#    existing CLICK userbtn / READTABLE table rows (old rows 7-8) down
#    to rows 8-9.
# ------------------------------------------------------------------
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "Read User Table"

# ------------------------------------------------------------------
# 2. Append the new "Add user" section (rows 10-20).
# ------------------------------------------------------------------
$ws.Range("A10").Value = "Add user"

$ws.Range("B11").Value = "CLICK"
$ws.Range("C11").Value = "addUser"
$ws.Range("D11").Value = "XPATH"

$ws.Range("B12").Value = "ENTERTEXT"
$ws.Range("C12").Value = "name"
$ws.Range("D12").Value = "ID"
$ws.Range("E12").Value = "mangesh"

$ws.Range("B13").Value = "ENTERTEXT"
$ws.Range("C13").Value = "mobile"
$ws.Range("D13").Value = "ID"
$ws.Range("E13").Value = 7894561230

$ws.Range("B14").Value = "ENTERTEXT"
$ws.Range("C14").Value = "email"
$ws.Range("D14").Value = "ID"
$ws.Range("E14").Value = "mangesh@gmail.com"

$ws.Range("B15").Value = "ENTERTEXT"
$ws.Range("C15").Value = "course"
$ws.Range("D15").Value = "ID"
$ws.Range("E15").Value = "selenium"

$ws.Range("B16").Value = "CLICK"
$ws.Range("C16").Value = "gender"
$ws.Range("D16").Value = "ID"

$ws.Range("B17").Value = "SELECTVALUE"
$ws.Range("C17").Value = "state"
$ws.Range("D17").Value = "XPATH"
$ws.Range("E17").Value = "HP"

$ws.Range("B18").Value = "ENTERTEXT"
$ws.Range("C18").Value = "pwd"
$ws.Range("D18").Value = "ID"
$ws.Range("E18").Value = 123456

$ws.Range("B19").Value = "CLICK"
$ws.Range("C19").Value = "submit"
$ws.Range("D19").Value = "ID"

$ws.Range("B20").Value = "ALERTHANDLE"

# ------------------------------------------------------------------
# 3. Normalize formatting: copy the plain bordered body-row look (from
#    an existing body row) onto every new cell, then re-apply the
#    centered alignment used across the whole table.
# ------------------------------------------------------------------
$ws.Range("A2:E2").Copy()
$ws.Range("A7:E20").PasteSpecial(-4122)

$body = $ws.Range("A1:E20")
$body.HorizontalAlignment = -4108
$body.VerticalAlignment = -4108

# ------------------------------------------------------------------
# 4. Sheet-level bookkeeping that mirrors the recorded edit: dimension
#    grows to A1:E20 automatically; move the active selection like the
#    author left it.
# ------------------------------------------------------------------
$ws.Range("E23").Select()
